$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 53: new bibliography entry ---
$ws.Range("A53").Value = "Evolution and behavioural responses to human-induced rapid environmental change"
$ws.Range("B53").Value = "Evolutionary application"
$ws.Range("C53").Value = 2011
$ws.Range("D53").Value = "Sih et al."

$ws.Range("F19").Copy() | Out-Null
$ws.Range("F53").PasteSpecial(-4122) | Out-Null
$ws.Range("F53").Value = "no"

$ws.Range("G53").Value = "Evolution in response to human-induced changes"
$ws.Range("I53").Value = "yes"
$ws.Range("J53").Value = "yes"
$ws.Range("K53").Value = "yes"

# --- Row 54: new bibliography entry ---
$ws.Range("A54").Value = "Ecological and evolutionary traps"
$ws.Range("B54").Value = "Trends Ecol. Evol."
$ws.Range("C54").Value = 2002
$ws.Range("D54").Value = "Schlaepfer et al."

$ws.Range("F19").Copy() | Out-Null
$ws.Range("F54").PasteSpecial(-4122) | Out-Null
$ws.Range("F54").Value = "no"

$ws.Range("G54").Value = "Evolutionary traps"
$ws.Range("I54").Value = "yes"
$ws.Range("J54").Value = "yes"

$ws.Range("K53").Copy() | Out-Null
$ws.Range("K54").PasteSpecial(-4122) | Out-Null
$ws.Range("K54").Value = "yes"

# --- Row 55: new bibliography entry ---
$ws.Range("A55").Value = "Ecological novelty and the emergence of evolutionary traps"
$ws.Range("B55").Value = "Trends Ecol. Evol."
$ws.Range("C55").Value = 2013
$ws.Range("D55").Value = "Robertson, Rehage & Sih"

$ws.Range("F19").Copy() | Out-Null
$ws.Range("F55").PasteSpecial(-4122) | Out-Null
$ws.Range("F55").Value = "no"

$ws.Range("K53").Copy() | Out-Null
$ws.Range("G55").PasteSpecial(-4122) | Out-Null
$ws.Range("G55").Value = "Evolutionary traps"

$ws.Range("I55").Value = "yes"
$ws.Range("J55").Value = "yes"

$ws.Range("K53").Copy() | Out-Null
$ws.Range("K55").PasteSpecial(-4122) | Out-Null
$ws.Range("K55").Value = "yes"

$excel.CutCopyMode = 0

# --- Refresh calculation so the summary totals (N19/N20/N22/N24) update ---
$excel.Calculate()

# --- Update the view state to match the final selection / scroll position ---
$excel.ActiveWindow.ScrollRow = 21
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("K55").Select() | Out-Null
